$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (currently holds the "Jelleg" values)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Move existing values (header + data) from column A to column B
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 1).Value2
}

# Set the new header for column A
$ws.Cells.Item(1, 1).Value2 = "ID"

# Fill column A with the sequential numeric id (1..N) for the data rows
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 1
}

# Restore the selection as recorded in the workbook after the edit
$ws.Range("F32").Select()
